# Apply the edit described by the diff:
#  - Sheet now only contains a single header row (row 1) spanning B1:K1
#  - Former data rows (2-5) are removed entirely
#  - Four new "Unnamed" index columns are inserted before the old columns
#  - The old "mdescription" / "IP" columns (G1/H1) are dropped
#  - Old A1 is no longer part of the used range

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Carry the existing header style (style index 1: bold, centered + top
# aligned, thin-bordered) over to the new header range B1:K1 *before* we
# touch any values, so the style table itself is left untouched (no new
# font/xf entries get minted) and the cells simply pick up the style that
# used to live on A1:H1.
$ws.Range("A1:H1").Copy() | Out-Null
$ws.Range("B1:K1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# Write the new header row text.
$ws.Range("B1").Value = "Unnamed: 0.3"
$ws.Range("C1").Value = "Unnamed: 0.2"
$ws.Range("D1").Value = "Unnamed: 0.1"
$ws.Range("E1").Value = "Unnamed: 0"
$ws.Range("F1").Value = "Machinetype"
$ws.Range("G1").Value = "Protocol"
$ws.Range("H1").Value = "Username"
$ws.Range("I1").Value = "Password"
$ws.Range("J1").Value = "Endpoint"
$ws.Range("K1").Value = "AddressNs"

# Drop what's no longer part of the sheet: the old A1 cell (now unused)
# and the old data rows 2-5 (Arburg/Fanuc rows + their formatting).
$ws.Range("A1").Clear() | Out-Null
$ws.Rows("2:5").Clear() | Out-Null
